$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Empresa1"

$ws.Range("A1").Value = "Exame"
$ws.Range("B1").Value = "Valor"

$ws.Range("A2").Value = "teste1"
$ws.Range("A3").Value = "teste2"
$ws.Range("B3").Value = 13
$ws.Range("A4").Value = "teste3"
$ws.Range("A5").Value = "teste4"
